{"js": "// TS Kramam 1 and 2 final - 16/10/2021\n//\n// The opening title paragraph currently reads:\n//   \"TS Krama Paatam \u2013 TS 1.4 Tamil Corrections \u2013 Observed till ?????\"\n// followed by a separate, centered paragraph:\n//   \"(ignore those which are already incorporated in your book's version and date)\"\n//\n// This edit:\n//   1. Fills in the placeholder \"?????\" date with \"31st Oct 2021\" (appended to\n//      the \"till \" run) and removes the red-highlighted placeholder run's\n//      leftover formatting (it becomes a plain bold space).\n//   2. Merges the title paragraph with the \"(ignore those...)\" paragraph by\n//      replacing the paragraph break with a manual line break (<w:br/>), and\n//      moves that second paragraph's center alignment up onto the (now single)\n//      merged paragraph.\n//   3. Strips the now-superfluous size/underline paragraph-mark formatting\n//      from the merged paragraph's mark run properties.\n\nconst body = context.document.body;\n\n// Locate the two paragraphs involved: the title paragraph (with the \"?????\"\n// placeholder) and the very next paragraph (\"(ignore those ...)\").\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst notePara = paragraphs.items[1];\n\nif (titlePara.text.indexOf(\"Observed till ?????\") === -1) {\n  throw new Error(\"Unexpected document content: title paragraph not found\");\n}\nif (!/ignore those which are already incorporated/.test(notePara.text)) {\n  throw new Error(\"Unexpected document content: note paragraph not found\");\n}\n\n// A range spanning from the very start of the title paragraph to the very\n// end of the note paragraph -- i.e. both paragraphs, including the\n// paragraph break between them.\nconst fullRange = titlePara.getRange(\"Start\").expandTo(notePara.getRange(\"End\"));\n\n// Rebuild that whole span as a single paragraph:\n//  - identical run formatting for the existing title text,\n//  - \"till \" extended to \"till 31st Oct 2021\",\n//  - the \"?????\" run collapsed to a single (unhighlighted, un-sized,\n//    un-underlined) bold space,\n//  - a manual line break in place of the old paragraph mark,\n//  - the note-paragraph text carried over unchanged,\n//  - center alignment (formerly on the note paragraph) now on this merged\n//    paragraph, with the paragraph-mark run properties trimmed to just\n//    bold/bold-complex-script.\nconst mergedParagraphOoxml = `\n  <w:p w14:paraId=\"775B95AA\" w14:textId=\"2B1F7482\" w:rsidR=\"009B3E2F\" w:rsidRPr=\"002F76B4\" w:rsidRDefault=\"009B3E2F\" w:rsidP=\"009B3E2F\">\n    <w:pPr>\n      <w:jc w:val=\"center\"/>\n      <w:rPr><w:b/><w:bCs/></w:rPr>\n    </w:pPr>\n    <w:r w:rsidRPr=\"002F76B4\">\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t xml:space=\"preserve\">TS </w:t>\n    </w:r>\n    <w:r>\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t>Krama</w:t>\n    </w:r>\n    <w:r w:rsidRPr=\"002F76B4\">\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t xml:space=\"preserve\"> Paatam &#8211; TS 1.</w:t>\n    </w:r>\n    <w:r>\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t>4</w:t>\n    </w:r>\n    <w:r w:rsidRPr=\"002F76B4\">\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t xml:space=\"preserve\"> </w:t>\n    </w:r>\n    <w:r>\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t>Tamil</w:t>\n    </w:r>\n    <w:r w:rsidRPr=\"002F76B4\">\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t xml:space=\"preserve\"> Corrections &#8211;</w:t>\n    </w:r>\n    <w:r>\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t xml:space=\"preserve\"> </w:t>\n    </w:r>\n    <w:r w:rsidRPr=\"002F76B4\">\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t xml:space=\"preserve\">Observed </w:t>\n    </w:r>\n    <w:r w:rsidRPr=\"002F76B4\">\n      <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n      <w:t>till 31st Oct 2021</w:t>\n    </w:r>\n    <w:r w:rsidRPr=\"00D814D0\">\n      <w:rPr><w:b/><w:bCs/></w:rPr>\n      <w:t xml:space=\"preserve\"> </w:t>\n    </w:r>\n    <w:r w:rsidRPr=\"00D814D0\">\n      <w:rPr><w:b/><w:bCs/></w:rPr>\n      <w:br/>\n    </w:r>\n    <w:r w:rsidRPr=\"002F76B4\">\n      <w:rPr><w:b/><w:bCs/></w:rPr>\n      <w:t>(ignore those which are already incorporated in your book&#8217;s version and date)</w:t>\n    </w:r>\n  </w:p>`;\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + mergedParagraphOoxml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nfullRange.insertOoxml(ooxmlPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# TS Kramam 1 and 2 final - 16/10/2021\n#\n# The opening title paragraph currently reads:\n#   \"TS Krama Paatam - TS 1.4 Tamil Corrections - Observed till ?????\"\n# followed by a separate, centered paragraph:\n#   \"(ignore those which are already incorporated in your book's version and date)\"\n#\n# This edit:\n#   1. Fills in the placeholder \"?????\" date with \"31st Oct 2021\" (appended to\n#      the \"till \" run) and removes the red-highlighted placeholder run's\n#      leftover formatting (it becomes a plain bold space).\n#   2. Merges the title paragraph with the \"(ignore those...)\" paragraph by\n#      replacing the paragraph break with a manual line break (<w:br/>), and\n#      moves that second paragraph's center alignment up onto the (now single)\n#      merged paragraph.\n#   3. Strips the now-superfluous size/underline paragraph-mark formatting\n#      from the merged paragraph's mark run properties.\n\n$d = $word.ActiveDocument\n\n$titlePara = $d.Paragraphs(1)\n$notePara = $d.Paragraphs(2)\n\nif ($titlePara.Range.Text -notmatch \"Observed till \\?\\?\\?\\?\\?\") {\n    throw \"Unexpected document content: title paragraph not found\"\n}\nif ($notePara.Range.Text -notmatch \"ignore those which are already incorporated\") {\n    throw \"Unexpected document content: note paragraph not found\"\n}\n\n# A range spanning from the very start of the title paragraph to the very\n# end of the note paragraph -- i.e. both paragraphs, including the\n# paragraph break between them.\n$fullRange = $d.Range($titlePara.Range.Start, $notePara.Range.End)\n\n# Rebuild that whole span as a single paragraph:\n#  - identical run formatting for the existing title text,\n#  - \"till \" extended to \"till 31st Oct 2021\",\n#  - the \"?????\" run collapsed to a single (unhighlighted, un-sized,\n#    un-underlined) bold space,\n#  - a manual line break in place of the old paragraph mark,\n#  - the note-paragraph text carried over unchanged,\n#  - center alignment (formerly on the note paragraph) now on this merged\n#    paragraph, with the paragraph-mark run properties trimmed to just\n#    bold/bold-complex-script.\n$mergedParagraphOoxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p w14:paraId=\"775B95AA\" w14:textId=\"2B1F7482\" w:rsidR=\"009B3E2F\" w:rsidRPr=\"002F76B4\" w:rsidRDefault=\"009B3E2F\" w:rsidP=\"009B3E2F\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n            <w:pPr>\n              <w:jc w:val=\"center\"/>\n              <w:rPr><w:b/><w:bCs/></w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"002F76B4\">\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t xml:space=\"preserve\">TS </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t>Krama</w:t>\n            </w:r>\n            <w:r w:rsidRPr=\"002F76B4\">\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t xml:space=\"preserve\"> Paatam &#8211; TS 1.</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t>4</w:t>\n            </w:r>\n            <w:r w:rsidRPr=\"002F76B4\">\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t>Tamil</w:t>\n            </w:r>\n            <w:r w:rsidRPr=\"002F76B4\">\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t xml:space=\"preserve\"> Corrections &#8211;</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:r w:rsidRPr=\"002F76B4\">\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t xml:space=\"preserve\">Observed </w:t>\n            </w:r>\n            <w:r w:rsidRPr=\"002F76B4\">\n              <w:rPr><w:b/><w:bCs/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:u w:val=\"single\"/></w:rPr>\n              <w:t>till 31st Oct 2021</w:t>\n            </w:r>\n            <w:r w:rsidRPr=\"00D814D0\">\n              <w:rPr><w:b/><w:bCs/></w:rPr>\n              <w:t xml:space=\"preserve\"> </w:t>\n            </w:r>\n            <w:r w:rsidRPr=\"00D814D0\">\n              <w:rPr><w:b/><w:bCs/></w:rPr>\n              <w:br/>\n            </w:r>\n            <w:r w:rsidRPr=\"002F76B4\">\n              <w:rPr><w:b/><w:bCs/></w:rPr>\n              <w:t>(ignore those which are already incorporated in your book&#8217;s version and date)</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$fullRange.InsertXML($mergedParagraphOoxml)\n"}
